$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 85

$ws.Cells.Item($row, 2).Value = "BETA"

$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "70000003601"
$ws.Cells.Item($row, 3).ClearFormats()

$ws.Cells.Item($row, 4).Value = "HHzlF+test@cyi.com"

$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Cells.Item($row, 5).Value = "1234567"
$ws.Cells.Item($row, 5).ClearFormats()

$ws.Cells.Item($row, 6).Value = "Australia"
$ws.Cells.Item($row, 7).Value = "88cf812c-5512-4c40-9a3d-95170336f46a"
